# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# 56f5903a-27d0-450b-b254-6d866f433341.md file is now "Ready for handoff"
# (instead of "Handed back: in sync with en-US"), with refreshed timestamps
# and a new error detail describing the stale handback version.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/19b692ac8810dc3fed990872e8e50d469dafaa2b/e2e/56f5903a-27d0-450b-b254-6d866f433341.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09b54e96613b586680d066ed53de9e8703f26c92/e2e/56f5903a-27d0-450b-b254-6d866f433341.md."

# --- Overview sheet: row 3 corresponds to 56f5903a-27d0-450b-b254-6d866f433341.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-31 00:51:46"

# --- zh-cn sheet: row 3 corresponds to 56f5903a-27d0-450b-b254-6d866f433341.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-31 00:51:41"
$wsZhCn.Range("P3").Value = $errorDetail
# 39.1667 is the ColumnWidth input that Excel's char-width rounding maps to a
# stored OOXML width of exactly 40 (matches the wider Error Detail column).
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1667

# --- de-de sheet: row 3 corresponds to 56f5903a-27d0-450b-b254-6d866f433341.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-31 00:51:46"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1667
